# Update odds data on "Sheet1" of the FlashScore weekly-games workbook
# (Jogos_da_Semana_FlashScore_2024-11-17.xlsx) to match the latest
# scraped values for rows 2, 3, 4, 11, 12, 13, 14, 15 and 16.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - San Lorenzo vs Racing Club
$ws.Range("O2").Value = 1.73
$ws.Range("P2").Value = 2

# Row 3 - Deportivo Madryn vs Gimnasia Mendoza
$ws.Range("G3").Value = 1.75
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 5.5
$ws.Range("J3").Value = 2.6
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53
$ws.Range("X3").Value = 6.5
$ws.Range("Z3").Value = 15
$ws.Range("AG3").Value = 10
$ws.Range("AI3").Value = 19
$ws.Range("AJ3").Value = 51
$ws.Range("AZ3").Value = 126

# Row 4 - San Martin S.J. vs All Boys
$ws.Range("M4").Value = 1.17
$ws.Range("N4").Value = 5
$ws.Range("O4").Value = 1.73
$ws.Range("P4").Value = 2
$ws.Range("S4").Value = 1.75
$ws.Range("T4").Value = 2.05

# Row 11 - General Caballero JLM vs Libertad Asuncion
$ws.Range("G11").Value = 2.8
$ws.Range("I11").Value = 2.4
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("AW11").Value = 4.33

# Row 12 - Tacuary vs Nacional Asuncion
$ws.Range("G12").Value = 4
$ws.Range("H12").Value = 3.3
$ws.Range("I12").Value = 1.95
$ws.Range("J12").Value = 4.5
$ws.Range("K12").Value = 2.05
$ws.Range("L12").Value = 2.63
$ws.Range("N12").Value = 9
$ws.Range("O12").Value = 1.36
$ws.Range("P12").Value = 3
$ws.Range("U12").Value = 1.91
$ws.Range("V12").Value = 1.8
$ws.Range("W12").Value = 10
$ws.Range("X12").Value = 19
$ws.Range("Z12").Value = 41
$ws.Range("AA12").Value = 34
$ws.Range("AH12").Value = 8.5
$ws.Range("AJ12").Value = 17
$ws.Range("AK12").Value = 17
$ws.Range("AN12").Value = 5.5
$ws.Range("AO12").Value = 23
$ws.Range("AQ12").Value = 81
$ws.Range("AR12").Value = 101
$ws.Range("AS12").Value = 251
$ws.Range("AX12").Value = 11
$ws.Range("AZ12").Value = 41
$ws.Range("BB12").Value = 201

# Row 13 - Zaragoza vs Malaga
$ws.Range("G13").Value = 1.85
$ws.Range("H13").Value = 3.1
$ws.Range("I13").Value = 4.75
$ws.Range("J13").Value = 2.63
$ws.Range("L13").Value = 5.5
$ws.Range("M13").Value = 1.11
$ws.Range("N13").Value = 6.5
$ws.Range("X13").Value = 7.5
$ws.Range("Z13").Value = 15
$ws.Range("AD13").Value = 6.5
$ws.Range("AE13").Value = 21
$ws.Range("AF13").Value = 81
$ws.Range("AG13").Value = 10
$ws.Range("AI13").Value = 17
$ws.Range("AN13").Value = 3.6
$ws.Range("AX13").Value = 29

# Row 14 - Racing Montevideo vs Danubio
$ws.Range("G14").Value = 2.2
$ws.Range("I14").Value = 3.7
$ws.Range("K14").Value = 1.83
$ws.Range("L14").Value = 4.75
$ws.Range("Z14").Value = 21
$ws.Range("AC14").Value = 5.5
$ws.Range("AO14").Value = 15
$ws.Range("AV14").Value = 81
$ws.Range("AX14").Value = 23
$ws.Range("AZ14").Value = 81

# Row 15 - Penarol vs Defensor Sp.
$ws.Range("G15").Value = 1.42
$ws.Range("I15").Value = 7
$ws.Range("K15").Value = 2.1
$ws.Range("L15").Value = 8.5
$ws.Range("M15").Value = 1.06
$ws.Range("N15").Value = 10
$ws.Range("AD15").Value = 9.5
$ws.Range("AG15").Value = 12
$ws.Range("AH15").Value = 34
$ws.Range("AI15").Value = 23
$ws.Range("AL15").Value = 67
$ws.Range("AP15").Value = 23
$ws.Range("AQ15").Value = 23
$ws.Range("AW15").Value = 9
$ws.Range("AX15").Value = 41
$ws.Range("AZ15").Value = 251
$ws.Range("BA15").Value = 301

# Row 16 - Estudiantes Merida vs Carabobo
$ws.Range("G16").Value = 3.35
$ws.Range("H16").Value = 3.25
$ws.Range("I16").Value = 2.07
$ws.Range("J16").Value = 3.8
$ws.Range("K16").Value = 2.07
$ws.Range("L16").Value = 2.65
$ws.Range("O16").Value = 1.27
$ws.Range("P16").Value = 3.1
$ws.Range("Q16").Value = 1.85
$ws.Range("R16").Value = 1.85
$ws.Range("W16").Value = 10.5
$ws.Range("X16").Value = 18.5
$ws.Range("Y16").Value = 11.5
$ws.Range("Z16").Value = 45
$ws.Range("AA16").Value = 29
$ws.Range("AB16").Value = 35
$ws.Range("AE16").Value = 13
$ws.Range("AG16").Value = 7.8
$ws.Range("AH16").Value = 10.25
$ws.Range("AI16").Value = 8.5
$ws.Range("AJ16").Value = 19.5
$ws.Range("AK16").Value = 16
$ws.Range("AL16").Value = 25
$ws.Range("AM16").Value = 400
$ws.Range("AN16").Value = 5.2
$ws.Range("AO16").Value = 18.5
$ws.Range("AP16").Value = 24
$ws.Range("AQ16").Value = 90
$ws.Range("AR16").Value = 120
$ws.Range("AU16").Value = 6.8
$ws.Range("AW16").Value = 3.95
$ws.Range("AX16").Value = 10.5
$ws.Range("AY16").Value = 18.5
$ws.Range("AZ16").Value = 40
$ws.Range("BA16").Value = 70
